$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 701.875
$ws.Range("I6").Value = 205
$ws.Range("K6").Value = 615
$ws.Range("M6").Value = -503
$ws.Range("H113").Value = 2205.15
$ws.Range("I113").Value = 2018.9375
$ws.Range("J113").Value = 2950
$ws.Range("K113").Value = 2018.9375
$ws.Range("L113").Value = 2950
$ws.Range("M113").Value = 1235.0625
$ws.Range("N113").Value = -9458
$ws.Range("H116").Value = 2014.1052
$ws.Range("I116").Value = 1783.3334
$ws.Range("J116").Value = 2120.6155
$ws.Range("K116").Value = 1783.3334
$ws.Range("L116").Value = 2120.6155
$ws.Range("M116").Value = 1658.6666
$ws.Range("N116").Value = -9004.6155
$ws.Range("H132").Value = 3190.15
$ws.Range("I132").Value = 2808.5557
$ws.Range("J132").Value = 6624.5
$ws.Range("K132").Value = 8425.667099999999
$ws.Range("L132").Value = 19873.5
$ws.Range("M132").Value = -5895.667099999999
$ws.Range("N132").Value = -24933.5
$ws.Range("H134").Value = 50120.453
$ws.Range("I134").Value = 20709
$ws.Range("J134").Value = 57179.2
$ws.Range("K134").Value = 20709
$ws.Range("L134").Value = 57179.2
$ws.Range("M134").Value = -15639
$ws.Range("N134").Value = -67319.2
$ws.Range("H137").Value = 10527157
$ws.Range("I137").Value = 898.3570999999999
$ws.Range("K137").Value = 2695.0713
$ws.Range("M137").Value = -145.0712999999996
$ws.Range("H139").Value = 36211.668
$ws.Range("J139").Value = 43963
$ws.Range("L139").Value = 43963
$ws.Range("N139").Value = -54243
$ws.Range("H140").Value = 48500
$ws.Range("J140").Value = 48500
$ws.Range("L140").Value = 48500
$ws.Range("N140").Value = -58860
$ws.Range("H141").Value = 1648.4375
$ws.Range("I141").Value = 1491.6666
$ws.Range("J141").Value = 4000
$ws.Range("K141").Value = 4474.9998
$ws.Range("L141").Value = 12000
$ws.Range("M141").Value = 705.0002000000004
$ws.Range("N141").Value = -22360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H29").Value = 2500
$ws.Range("I29").Value = 2500
$ws.Range("K29").Value = 2500
$ws.Range("M29").Value = -2192
$ws.Range("H61").Value = 21741466
$ws.Range("I61").Value = 25002506
$ws.Range("J61").Value = 1198
$ws.Range("K61").Value = 25002506
$ws.Range("L61").Value = 1198
$ws.Range("M61").Value = -25002294
$ws.Range("N61").Value = -1622
$ws.Range("H74").Value = 14709667
$ws.Range("I74").Value = 20835432
$ws.Range("J74").Value = 7832.8
$ws.Range("K74").Value = 20835432
$ws.Range("L74").Value = 7832.8
$ws.Range("M74").Value = -20834558
$ws.Range("N74").Value = -9580.799999999999
$ws.Range("H77").Value = 14709667
$ws.Range("I77").Value = 20835432
$ws.Range("J77").Value = 7832.8
$ws.Range("K77").Value = 104177160
$ws.Range("L77").Value = 39164
$ws.Range("M77").Value = -104172792
$ws.Range("N77").Value = -47900
$ws.Range("H132").Value = 7815219
$ws.Range("I132").Value = 13891130
$ws.Range("J132").Value = 3333.4285
$ws.Range("K132").Value = 41673390
$ws.Range("L132").Value = 10000.2855
$ws.Range("M132").Value = -41670860
$ws.Range("N132").Value = -15060.2855
$ws.Range("H136").Value = 21741466
$ws.Range("I136").Value = 25002506
$ws.Range("J136").Value = 1198
$ws.Range("K136").Value = 75007518
$ws.Range("L136").Value = 3594
$ws.Range("M136").Value = -75004968
$ws.Range("N136").Value = -8694

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H116").Value = 50700
$ws.Range("J116").Value = 50700
$ws.Range("L116").Value = 50700
$ws.Range("N116").Value = -59878
$ws.Range("H132").Value = 53325.9
$ws.Range("J132").Value = 53325.9
$ws.Range("L132").Value = 53325.9
$ws.Range("N132").Value = -63445.9
$ws.Range("H134").Value = 2633.25
$ws.Range("I134").Value = 2001.1364
$ws.Range("J134").Value = 6109.875
$ws.Range("K134").Value = 6003.4092
$ws.Range("L134").Value = 18329.625
$ws.Range("M134").Value = -3468.4092
$ws.Range("N134").Value = -23399.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1665.7632
$ws.Range("I58").Value = 672.53845
$ws.Range("J58").Value = 3817.75
$ws.Range("K58").Value = 672.53845
$ws.Range("L58").Value = 3817.75
$ws.Range("M58").Value = -469.53845
$ws.Range("N58").Value = -4223.75
$ws.Range("H70").Value = 39836.668
$ws.Range("J70").Value = 39836.668
$ws.Range("L70").Value = 39836.668
$ws.Range("N70").Value = -40466.668
$ws.Range("H73").Value = 39836.668
$ws.Range("J73").Value = 39836.668
$ws.Range("L73").Value = 39836.668
$ws.Range("N73").Value = -42020.668
$ws.Range("H99").Value = 1490.1765
$ws.Range("I99").Value = 1343.3
$ws.Range("J99").Value = 1700
$ws.Range("K99").Value = 1343.3
$ws.Range("L99").Value = 1700
$ws.Range("M99").Value = 154.7
$ws.Range("N99").Value = -4696
$ws.Range("H126").Value = 1490.1765
$ws.Range("I126").Value = 1343.3
$ws.Range("J126").Value = 1700
$ws.Range("K126").Value = 4029.9
$ws.Range("L126").Value = 5100
$ws.Range("M126").Value = -1559.9
$ws.Range("N126").Value = -10040
$ws.Range("H134").Value = 2260.0588
$ws.Range("I134").Value = 2215.923
$ws.Range("J134").Value = 2403.5
$ws.Range("K134").Value = 6647.768999999999
$ws.Range("L134").Value = 7210.5
$ws.Range("M134").Value = -4112.768999999999
$ws.Range("N134").Value = -12280.5
$ws.Range("H136").Value = 1665.7632
$ws.Range("I136").Value = 672.53845
$ws.Range("J136").Value = 3817.75
$ws.Range("K136").Value = 2017.61535
$ws.Range("L136").Value = 11453.25
$ws.Range("M136").Value = 532.38465
$ws.Range("N136").Value = -16553.25
$ws.Range("H140").Value = 39603.934
$ws.Range("J140").Value = 39603.934
$ws.Range("L140").Value = 39603.934
$ws.Range("N140").Value = -49963.934

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 933
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 933
$ws.Range("K68").Value = 0
$ws.Range("M68").Value = 2799
$ws.Range("N68").Value = -4421
$ws.Range("H71").Value = 933
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 933
$ws.Range("K71").Value = 0
$ws.Range("M71").Value = 8397
$ws.Range("N71").Value = -16509

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4512.5186
$ws.Range("I126").Value = 3261.75
$ws.Range("J126").Value = 5513.1333
$ws.Range("K126").Value = 9785.25
$ws.Range("L126").Value = 16539.3999
$ws.Range("M126").Value = -7315.25
$ws.Range("N126").Value = -21479.3999
$ws.Range("H132").Value = 5988.091
$ws.Range("I132").Value = 5012.25
$ws.Range("J132").Value = 6545.7144
$ws.Range("K132").Value = 15036.75
$ws.Range("L132").Value = 19637.1432
$ws.Range("M132").Value = -12506.75
$ws.Range("N132").Value = -24697.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5006.9395
$ws.Range("I7").Value = 5508.857
$ws.Range("J7").Value = 4637.1055
$ws.Range("K7").Value = 5508.857
$ws.Range("L7").Value = 4637.1055
$ws.Range("M7").Value = -5396.857
$ws.Range("N7").Value = -4861.1055
$ws.Range("H16").Value = 1808.0769
$ws.Range("I16").Value = 2130.3
$ws.Range("J16").Value = 734
$ws.Range("K16").Value = 2130.3
$ws.Range("L16").Value = 734
$ws.Range("M16").Value = -1960.3
$ws.Range("N16").Value = -1074
$ws.Range("H40").Value = 5634
$ws.Range("I40").Value = 5713.923
$ws.Range("J40").Value = 5485.5713
$ws.Range("K40").Value = 5713.923
$ws.Range("L40").Value = 5485.5713
$ws.Range("M40").Value = -5577.923
$ws.Range("N40").Value = -5757.5713
$ws.Range("H55").Value = 564.1539
$ws.Range("I55").Value = 391.5
$ws.Range("J55").Value = 840.4
$ws.Range("K55").Value = 391.5
$ws.Range("L55").Value = 840.4
$ws.Range("M55").Value = -218.5
$ws.Range("N55").Value = -1186.4
$ws.Range("H93").Value = 998
$ws.Range("I93").Value = 847.5
$ws.Range("K93").Value = 847.5
$ws.Range("M93").Value = 400.5
$ws.Range("H126").Value = 5006.9395
$ws.Range("I126").Value = 5508.857
$ws.Range("J126").Value = 4637.1055
$ws.Range("K126").Value = 16526.571
$ws.Range("L126").Value = 13911.3165
$ws.Range("M126").Value = -14056.571
$ws.Range("N126").Value = -18851.3165
$ws.Range("H132").Value = 12580.24
$ws.Range("I132").Value = 10133.111
$ws.Range("J132").Value = 13956.75
$ws.Range("K132").Value = 30399.333
$ws.Range("L132").Value = 41870.25
$ws.Range("M132").Value = -27869.333
$ws.Range("N132").Value = -46930.25
$ws.Range("H139").Value = 59868
$ws.Range("J139").Value = 59868
$ws.Range("L139").Value = 59868
$ws.Range("N139").Value = -70148

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 23501.666
$ws.Range("I7").Value = 250
$ws.Range("K7").Value = 250
$ws.Range("M7").Value = -137
$ws.Range("H96").Value = 1412.1428
$ws.Range("I96").Value = 1094.1578
$ws.Range("J96").Value = 2083.4443
$ws.Range("K96").Value = 1094.1578
$ws.Range("L96").Value = 2083.4443
$ws.Range("M96").Value = 278.8422
$ws.Range("N96").Value = -4829.4443
$ws.Range("H123").Value = 25113.777
$ws.Range("J123").Value = 25113.777
$ws.Range("L123").Value = 25113.777
$ws.Range("N123").Value = -34913.777
